$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.137.63"
$ws.Range("E2").Value = "  +0.91%  "

$ws.Range("D3").Value = "3.021.08"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.48"
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.08"
$ws.Range("E6").Value = "  -3.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  -3.06%  "

$ws.Range("D9").Value = "3.029.11"
$ws.Range("E9").Value = "  -2.08%  "

$ws.Range("E10").Value = "  -1.28%  "

$ws.Range("E11").Value = "  -4.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.367"
$ws.Range("E12").Value = "  -2.05%  "

$ws.Range("D13").Value = "3.545.83"
$ws.Range("E13").Value = "  -2.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.124"
$ws.Range("E14").Value = "  -2.86%  "

$ws.Range("D15").Value = "63.226.75"
$ws.Range("E15").Value = "  +0.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "24.14"
$ws.Range("E16").Value = "  -1.32%  "

$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("D18").Value = "3.030.73"
$ws.Range("E18").Value = "  -2.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "397.56"
$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("E20").Value = "  +0.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.05"
$ws.Range("E21").Value = "  -2.02%  "

$ws.Range("E22").Value = "  -4.64%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.35"
$ws.Range("E24").Value = "  -2.93%  "

$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.191"
$ws.Range("E25").Value = "  -2.16%  "

$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.467"
$ws.Range("E26").Value = "  -1.76%  "

$ws.Range("D27").Value = "0.0₃0988"
$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("E28").Value = "  +1.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  -0.49%  "

$ws.Range("E32").Value = "  -1.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "160.20"
$ws.Range("E33").Value = "  +4.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("E34").Value = "  -0.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.11"
$ws.Range("E35").Value = "  +2.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.07"
$ws.Range("E36").Value = "  -1.57%  "

$ws.Range("E37").Value = "  +0.32%  "

$ws.Range("D38").Value = "2.546.68"
$ws.Range("E38").Value = "  -5.63%  "

$ws.Range("E39").Value = "  -3.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.96"
$ws.Range("E40").Value = "  -1.07%  "

$ws.Range("E41").Value = "  -0.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.66"
$ws.Range("E42").Value = "  -1.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.670"
$ws.Range("E43").Value = "  -2.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0602"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("E45").Value = "  -0.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.09"
$ws.Range("E46").Value = "  -2.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  -0.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.30"
$ws.Range("E48").Value = "  -1.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "269.78"
$ws.Range("E49").Value = "  -3.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0947"
$ws.Range("E50").Value = "  -2.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.49"
$ws.Range("E51").Value = "  +0.45%  "
